# Add a new product row ("ISIS HS 750g LEMON LIMITLESS") right above the
# "Le Chat HS 300 gr" row, and remove the two obsolete
# "ISIS HS POWDER SDM ..." rows further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5 (old row 5 "Le Chat HS 300 gr" and
# everything below it shifts down by one, formatting intact).
$ws.Rows(5).Insert()

# Give the new row the same look (borders) as the surrounding data rows by
# copying the formatting from row 4, then fill in its own values.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Range("A5").Value = 2875892
$ws.Range("B5").Value = "ISIS HS 750g LEMON LIMITLESS"

# Remove the two "ISIS HS POWDER SDM" rows (now at rows 11 and 12 after the
# insert above), shifting the rows below them up.
$ws.Range("A11:D12").Delete(-4162)

# Match the saved selection state.
$ws.Range("C13").Select()
